$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VA.MHV.bloodOxygenSat: narrow the multi-code "Code" list down to a single code
$ws.Range("E2").Value = "null#2708-6"

# VA.MHV.bloodPressure: convert to FHIR Core vital-signs profile conventions
$ws.Range("E3").Value = "null#85354-9"
$ws.Range("F3").Value = "http://hl7.org/fhir/ValueSet/observation-vitalsignresult (extensible)"
$ws.Range("G3").Value = "dateTimeĵ"
$ws.Range("H3").Value = "Quantityĵ"

# Blood pressure component rows (systolic/diastolic) move their LOINC codes
# from the "Code" column into the "Method" column, and pick up the ĵ-flagged
# Value Types used by the other vital-signs profiles.
$ws.Range("E4").Value = ""
$ws.Range("H4").Value = "Quantityĵ"
$ws.Range("K4").Value = "LOINC#8480-6"

$ws.Range("E5").Value = ""
$ws.Range("H5").Value = "Quantityĵ"
$ws.Range("K5").Value = "LOINC#8462-4"

# VA.MHV.bodyTemperature: code system switches from LOINC# to null#, and adopts
# the vitalsignresult value set / ĵ-flagged types
$ws.Range("E7").Value = "null#8310-5"
$ws.Range("F7").Value = "http://hl7.org/fhir/ValueSet/observation-vitalsignresult (extensible)"

# VA.MHV.bodyWeight: same treatment
$ws.Range("E8").Value = "null#29463-7"
$ws.Range("F8").Value = "http://hl7.org/fhir/ValueSet/observation-vitalsignresult (extensible)"
